# DOMA-11339: add "Archive date" column (Q) to the propertyMeter import example sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style + column width) from the existing last data column (P)
# onto the new column (Q), so the new column visually/structurally matches its
# neighbours (same borders/fill/number format, and the same 23.5-char width).
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("P2:P11").Copy()
$ws.Range("Q2:Q11").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("Q1:Q1").ColumnWidth = $ws.Range("P1").ColumnWidth

# Header for the new column
$ws.Range("Q1").Value = "Archive date"

# Example archive date value for the first data row only (rest stay blank,
# same as the other optional date columns in this template)
$ws.Range("Q2").Value = "2022-01-25"
